# feat: add 2022-Q1 data
#
# 1) Duplicate the "2021-Q4" sheet (same column layout) right after itself,
#    rename the copy to "2022-Q1", and overwrite its data with the new
#    quarter's fund breakdown.
# 2) Update the "总计" (summary) sheet: insert a new first data row for
#    2022-Q1 and push the existing rows down, refreshing the running index
#    column.

$wb = $excel.ActiveWorkbook

# --- Step 1: create the "2022-Q1" sheet from a copy of "2021-Q4" ---------
$wsQ4 = $wb.Worksheets.Item("2021-Q4")
$wsQ4.Copy($null, $wsQ4)
$wsNew = $wb.Worksheets.Item(3)
$wsNew.Name = "2022-Q1"

# Columns B-G hold text values (fund code / name / amounts as strings) in
# every quarter sheet, so keep them text here too - otherwise numeric-looking
# strings like "009686" or "0.0230" would be silently coerced to numbers and
# lose their leading/trailing zeros.
$wsNew.Range("B2:G3").NumberFormat = "@"

# Row 2
$wsNew.Cells.Item(2, 1).Value = 0
$wsNew.Cells.Item(2, 2).Value = "009686"
$wsNew.Cells.Item(2, 3).Value = "华夏磐利一年定期开放混合A"
$wsNew.Cells.Item(2, 4).Value = "16.02"
$wsNew.Cells.Item(2, 5).Value = "93.69"
$wsNew.Cells.Item(2, 6).Value = "4.51"
$wsNew.Cells.Item(2, 7).Value = "0.7225"
$wsNew.Cells.Item(2, 8).Value = 2

# Row 3
$wsNew.Cells.Item(3, 1).Value = 1
$wsNew.Cells.Item(3, 2).Value = "009687"
$wsNew.Cells.Item(3, 3).Value = "华夏磐利一年定期开放混合C"
$wsNew.Cells.Item(3, 4).Value = "0.51"
$wsNew.Cells.Item(3, 5).Value = "93.69"
$wsNew.Cells.Item(3, 6).Value = "4.51"
$wsNew.Cells.Item(3, 7).Value = "0.0230"
$wsNew.Cells.Item(3, 8).Value = 2

# --- Step 2: update the "总计" summary sheet ------------------------------
# Push the two existing rows down one slot (write bottom-up so nothing is
# clobbered before it is read) and add the new 2022-Q1 row on top, rather
# than using Rows.Insert() - Insert() drags the header row's bold/bordered
# format onto the freshly inserted row, which the source data rows don't
# have.
$wsTotal = $wb.Worksheets.Item("总计")

$wsTotal.Cells.Item(4, 1).Value = 2
$wsTotal.Cells.Item(4, 2).Value = "2021-Q3"
$wsTotal.Cells.Item(4, 3).Value = 2
$wsTotal.Cells.Item(4, 4).Value = 0.02
# Row 4 is brand new - copy row 3's formatting onto column A so the running
# index keeps the same (centered/bordered) look as the other index cells.
$wsTotal.Cells.Item(3, 1).Copy()
$wsTotal.Cells.Item(4, 1).PasteSpecial(-4122)

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2021-Q4"
$wsTotal.Cells.Item(3, 3).Value = 2
$wsTotal.Cells.Item(3, 4).Value = 0.27

$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(2, 2).Value = "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 2
$wsTotal.Cells.Item(2, 4).Value = 0.75

# Worksheet.Copy() makes the new copy the active tab as a side effect;
# restore the original active sheet so it is the only one left selected.
$wb.Worksheets.Item("2021-Q3").Activate()
